$d = $word.ActiveDocument

function Replace-Text([string]$old, [string]$new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2) | Out-Null
}

function Set-Bold([string]$text, [int]$bold) {
    $rng = $d.Content
    $found = $rng.Find.Execute($text, $true, $false, $false, $false, $false,
                                $true, 1, $false, "", 0)
    if ($found) {
        $rng.Font.Bold = $bold
    }
}

# 1. Product name / title
Replace-Text ": Mystisch Gewürz Premium Chai Tee" ": Mystic Spice Premium Chai Tee"

# 2. "Wichtige Merkmale:" -> "Wichtige Features:" (also becomes bold)
Replace-Text "Wichtige Merkmale:" "Wichtige Features:"
Set-Bold "Wichtige Features:" 1

# 3. "Authentic Blend" -> "Authentische Mischung"
Replace-Text "Authentic Blend" "Authentische Mischung"

# 4. Authentic blend description text
Replace-Text ": Unsere Chai ist eine harmonische Mischung aus Premium-Schwarzen Teeblättern und einer charakteristischen Auswahl an gemahlenen Gewürzen, darunter Knoblauch, Karamom, Gerinnsel, Ingwer und schwarzer Pfeffer." ": Unser Chai ist eine harmonische Mischung aus hochwertigen schwarzen Teeblättern und einer charakteristischen Auswahl an gemahlenen Gewürzen wie Zimt, Kardamom, Nelken, Ingwer und schwarzem Pfeffer."

# 5. "Gesundheitsfördernde Inhaltsstoffe: Jeder Bestandteil" -> "Gesundheitsfördernde Inhaltsstoffe"
Replace-Text "Gesundheitsfördernde Inhaltsstoffe: Jeder Bestandteil" "Gesundheitsfördernde Inhaltsstoffe"

# 6. Health ingredients description
Replace-Text " von Mystisch Gewürz-Chai-Tee wird für seine natürlichen Gesundheitlichen Vorteile ausgewählt." ": Alle Inhaltsstoffe des Mystic Spice Chai Tea werden aufgrund ihrer natürlichen gesundheitsfördernden Eigenschaften ausgewählt."

# 7. Rich aroma/taste description
Replace-Text ": Das warme, würzige Aroma und tiefe, belebende Geschmack unserer Chai machen es zum perfekten Getränk, um Ihren Tag zu beginnen oder sich am Abend zu entspannen." ": Das warme, würzige Aroma und tiefe, belebende Geschmack unseres Chai machen ihn zum perfekten Getränk, um in den Tag zu starten oder am Abend zu entspannen."

# 8. "Vielseitige Brauoptionen" -> "Vielfältige Zubereitungsmöglichkeiten"
Replace-Text "Vielseitige Brauoptionen" "Vielfältige Zubereitungsmöglichkeiten"

# 9. Brewing options description
Replace-Text ": Ob Sie Ihre Chai heiß dampfen, als erfrischender Eistee oder als cremefarbene Latte lieben, ist unsere Mischung vielseitig genug für jede Vorliebe." ": Ob Sie Ihren Chai dampfend heiß, als erfrischenden Eistee oder als cremigen Latte mögen – unsere Mischung ist vielseitig genug, um allen Vorlieben gerecht zu werden."

# 10. Sustainability description
Replace-Text ": Wir engagieren uns für Nachhaltigkeit, wir beziehen unsere Zutaten aus kleinflächigen Farmen, die ökologische Landwirtschaft betreiben, und sorgen nicht nur für die feinste Qualität, sondern auch für das Wohlergehen unseres Planeten." ": Da wir uns der Nachhaltigkeit verpflichtet haben, beziehen wir unsere Zutaten von kleinen Bauernhöfen, die ökologische Landwirtschaft betreiben. So garantieren wir nicht nur beste Qualität, sondern tragen auch zum Wohlergehen unseres Planeten bei."

# 11. Packaging description
Replace-Text ": Mystisch Gewürz-Chai-Tee kommt in wunderschön gestalteten, umweltfreundlichen Verpackungen, sodass es ein ideales Geschenk für Teeliebhaber oder ein luxuriöser Genuss für sich selbst ist." ": Mystic Spice Chai Tea wird in einer wunderschönen, umweltfreundlichen Verpackung geliefert, die ihn zu einem idealen Geschenk für Teeliebhaber oder zu einem luxuriösen Genuss für Sie selbst macht."

# 12. "Ideal für" description
Replace-Text ": Tee-Enthusiasten, gesundheitsbewusste Einzelpersonen, Liebhaber warmer, würziger Getränke und jeder, der die reichen Aromen der traditionellen indischen Chai erkunden möchte." ": Teeliebhaber, gesundheitsbewusste Menschen, Liebhaber von warmen, würzigen Getränken und alle, die den reichen Geschmack des traditionellen indischen Chai entdecken möchten."
